$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top for the column headers (A, B, C)
$ws.Rows.Item(1).Insert()

# Insert a new row before the current row 3 (which holds 11,12,13)
# for the section header (FIRST, SECOND, THIRD)
$ws.Rows.Item(4).Insert()

# Fill in the header rows
$ws.Range("A1").Value = "A"
$ws.Range("B1").Value = "B"
$ws.Range("C1").Value = "C"

$ws.Range("A4").Value = "FIRST"
$ws.Range("B4").Value = "SECOND"
$ws.Range("C4").Value = "THIRD"
